$d = $word.ActiveDocument

$replacements = @(
    @("955×7=6685", "217×8=1736"),
    @("947×8=7576", "741×8=5928"),
    @("704×9=6336", "682×8=5456"),
    @("527×5=2635", "485×9=4365"),
    @("518×2=1036", "859×3=2577"),
    @("649×3=1947", "343×8=2744"),
    @("308×7=2156", "621×3=1863"),
    @("513×2=1026", "569×4=2276"),
    @("276×9=2484", "781×8=6248"),
    @("439×3=1317", "582×8=4656"),
    @("874×6=5244", "193×8=1544"),
    @("889×6=5334", "277×5=1385"),
    @("608×6=3648", "872×6=5232"),
    @("690×6=4140", "419×2=838"),
    @("887×7=6209", "495×6=2970"),
    @("998×4=3992", "158×2=316"),
    @("755×2=1510", "152×7=1064"),
    @("792×7=5544", "473×9=4257"),
    @("310×3=930", "887×2=1774"),
    @("270×2=540", "400×5=2000"),
    @("453×2=906", "586×3=1758"),
    @("377×6=2262", "645×2=1290"),
    @("348×3=1044", "584×3=1752"),
    @("305×5=1525", "813×9=7317"),
    @("687×2=1374", "106×7=742")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
